$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Update the F-column (time_taken) timestamps on the "data" sheet ---
$ws1.Range("F2").Value = "2021-10-05 14:33:50.719924"
$ws1.Range("F3").Value = "2021-10-05 14:33:50.719935"
$ws1.Range("F4").Value = "2021-10-05 14:33:50.719939"
$ws1.Range("F5").Value = "2021-10-05 14:33:50.719941"
$ws1.Range("F6").Value = "2021-10-05 14:33:50.719944"
$ws1.Range("F7").Value = "2021-10-05 14:33:50.719947"
$ws1.Range("F8").Value = "2021-10-05 14:33:50.719950"
$ws1.Range("F9").Value = "2021-10-05 14:33:50.719952"
$ws1.Range("F10").Value = "2021-10-05 14:33:50.719955"
$ws1.Range("F11").Value = "2021-10-05 14:33:50.719958"
$ws1.Range("F12").Value = "2021-10-05 14:33:50.719960"
$ws1.Range("F13").Value = "2021-10-05 14:33:50.719963"

# --- 2. Add the new "metadata" sheet, positioned after "data" ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws2.Name = "metadata"

# Match the page margins used on the "data" sheet (PageSetup margins are in points).
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

# Header row (bold/border/center-top style matching the "data" sheet headers)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$ws2.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B2").Value = "Familial hypercholesterolaemia"
$ws2.Range("C2").Value = 333

# "0.18" must stay text (not become the number 0.18) but also must not pick
# up any style (target cell D2 has no "s" attribute) - force text via a
# temporary number format, then paste an unstyled donor cell's format over
# it to strip the format back off while keeping the text value intact.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.18"
$ws1.Range("E2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("E2").Value = "2021-09-04T06:00:55.517891Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:50.716058"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/333/?format=json"

$ws2.Range("A1").Select() | Out-Null
